$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (before Swampstriker) for Sellemental
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "Sellemental"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = "Minion"
$ws.Range("E9").Value = "Elemental"
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 17
$ws.Range("I9").Value = 0

# Update use_flg values for rows 2-4
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("H4").Value = "BattlecryMechanic"
$ws.Range("I4").Value = 1

$ws.Range("J1").Font.Name = "Menlo"
$ws.Range("J1").Font.Size = 14
$ws.Range("J1").Font.Color = 2562065

Write-Host "done"
